$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title / timestamp string (A1) ---
$ws.Cells.Item(1,1).Value = "Datos actualizados a 21 de Abril de 2020 a las 10:52"

# --- Update Suiza row (row 18) ---
$ws.Cells.Item(18,5).Value = 7908
$ws.Cells.Item(18,7).Value = 7
$ws.Cells.Item(18,8).Value = 1436

# --- Move Banglades row from 57 to 53 (new data), shifting Egipto/Sudafrica/Marruecos/Argentina down one row ---
# Row 53: Egipto -> Banglades (new numbers)
$ws.Cells.Item(53,1).Value = "Banglades"
$ws.Cells.Item(53,2).Value = 3382
$ws.Cells.Item(53,3).Value = 434
$ws.Cells.Item(53,4).Value = 87
$ws.Cells.Item(53,5).Value = 3185
$ws.Cells.Item(53,6).Value = 1
$ws.Cells.Item(53,7).Value = 9
$ws.Cells.Item(53,8).Value = 110

# Row 54: Sudafrica -> Egipto (old Egipto numbers)
$ws.Cells.Item(54,1).Value = "Egipto"
$ws.Cells.Item(54,2).Value = 3333
$ws.Cells.Item(54,3).Value = 0
$ws.Cells.Item(54,4).Value = 821
$ws.Cells.Item(54,5).Value = 2262
$ws.Cells.Item(54,6).Value = 0
$ws.Cells.Item(54,7).Value = 0
$ws.Cells.Item(54,8).Value = 250

# Row 55: Marruecos -> Sudafrica (old Sudafrica numbers)
$ws.Cells.Item(55,1).Value = "Sudafrica"
$ws.Cells.Item(55,2).Value = 3300
$ws.Cells.Item(55,3).Value = 0
$ws.Cells.Item(55,4).Value = 1055
$ws.Cells.Item(55,5).Value = 2187
$ws.Cells.Item(55,6).Value = 36
$ws.Cells.Item(55,7).Value = 0
$ws.Cells.Item(55,8).Value = 58

# Row 56: Argentina -> Marruecos (old Marruecos numbers)
$ws.Cells.Item(56,1).Value = "Marruecos"
$ws.Cells.Item(56,2).Value = 3046
$ws.Cells.Item(56,3).Value = 0
$ws.Cells.Item(56,4).Value = 350
$ws.Cells.Item(56,5).Value = 2553
$ws.Cells.Item(56,6).Value = 1
$ws.Cells.Item(56,7).Value = 0
$ws.Cells.Item(56,8).Value = 143

# Row 57: Banglades -> Argentina (old Argentina numbers)
$ws.Cells.Item(57,1).Value = "Argentina"
$ws.Cells.Item(57,2).Value = 3031
$ws.Cells.Item(57,3).Value = 0
$ws.Cells.Item(57,4).Value = 737
$ws.Cells.Item(57,5).Value = 2152
$ws.Cells.Item(57,6).Value = 123
$ws.Cells.Item(57,7).Value = 0
$ws.Cells.Item(57,8).Value = 142

# Row 58 (Tailandia) remains unchanged.
